$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# --- 1. New "Tube Suspension" calculation block (rows 22-27) ---
# Written first so its strings occupy shared-string indices 12-17.
$ws.Range("F22").Value = "Force"
$ws.Range("H22").Value = 800

$ws.Range("F23").Value = "Largeur"
$ws.Range("H23").Value = 0.227

$ws.Range("F24").Value = "Couple"
$ws.Range("H24").Formula = "=H22*H23"

$ws.Range("F25").Value = "Deplacement en Z"
$ws.Range("H25").Value = 0.003

$ws.Range("F26").Value = "Theta"
$ws.Range("H26").Formula = "=H25/H23"

$ws.Range("F27").Value = "K"
$ws.Range("H27").Formula = "=H24/H26"

# --- 2. Cell_Arr table: Front/Rear + A-Arm designations (rows 6-13) ---
# Column B (Front/Rear) is populated fully before column C so that
# "Front" / "Rear" occupy indices 18-19 ahead of the arm names (20-23).
$ws.Range("B6").Value = "Front"
$ws.Range("B7").Value = "Front"
$ws.Range("B8").Value = "Front"
$ws.Range("B9").Value = "Front"
$ws.Range("B10").Value = "Rear"
$ws.Range("B11").Value = "Rear"
$ws.Range("B12").Value = "Rear"
$ws.Range("B13").Value = "Rear"

$ws.Range("C6").Value = "Upper A-Arm Front"
$ws.Range("C7").Value = "Upper A-Arm Rear"
$ws.Range("C8").Value = "Lower A-Arm Front"
$ws.Range("C9").Value = "Lower A-Arm Rear"
$ws.Range("C10").Value = "Upper A-Arm Front"
$ws.Range("C11").Value = "Upper A-Arm Rear"
$ws.Range("C12").Value = "Lower A-Arm Front"
$ws.Range("C13").Value = "Lower A-Arm Rear"

# --- 3. Insert a "Deplacement en X" column in the header row (row 5) ---
# Shift the old G5/H5 header text one column to the right, then put the
# new header in G5 (copying F5's formatting), added to shared strings last.
$depY = $ws.Range("G5").Value()
$depZ = $ws.Range("H5").Value()
$ws.Range("H5").Value = $depY
$ws.Range("I5").Value = $depZ

$ws.Range("F5").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G5").Value = "Déplacement en X"

# --- 4. Extend the merged title block with a matching empty cell (col G) ---
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)

# --- 5. Column width adjustments to fit the new content ---
$ws.Columns("C").ColumnWidth = 15.665
$ws.Columns("G").ColumnWidth = 25.665
$ws.Columns("H").ColumnWidth = 14.83
$ws.Columns("I").ColumnWidth = 14.665

# --- 6. Restore selection to match the author's final cursor position ---
$excel.CutCopyMode = $false
$ws.Range("I20").Select()
